$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EV Annual")

# Re-add the per-model yearly-total formulas in column J (rows 2-29):
# total = sum of the 2011-2017 columns (C:I) for that row.
for ($r = 2; $r -le 29; $r++) {
    $ws.Range("J$r").Formula = "=SUM(C$r`:I$r)"
}

# Re-add the column totals in row 30 (C:I) summing each year column
# across all the model rows (2-29).
for ($col = "C"; $col -le "I"; $col = [char]([int][char]$col + 1)) {
    $ws.Range("$col`30").Formula = "=SUM($col`2:$col`29)"
}

# Grand total in J30 = sum of the yearly column totals.
$ws.Range("J30").Formula = "=SUM(C30:I30)"

# A blank, number-formatted cell was left behind at L22 (stray selection
# left over from re-entering the data); also move the active selection
# there to match.
$ws.Range("L22").NumberFormat = "#,##0"
[void]$ws.Range("L22").Select()
